$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 7293
$ws.Range("F4").Value = 7384
$ws.Range("G4").Value = 7532
$ws.Range("H4").Value = 7652
$ws.Range("I4").Value = 7801
$ws.Range("J4").Value = 8086
$ws.Range("K4").Value = 8158

$ws.Range("E4:K4").Select()
